$wb = $excel.ActiveWorkbook

$ws3 = $wb.Worksheets.Item("Sheet3")
$ws3.Range("B1").Value = "/a/c/2::Vector"

$ws2 = $wb.Worksheets.Item("Sheet2")
$ws2.Range("B1").Value = "/a/c/1::Vector"

$ws3.PageSetup.PaperSize = 9
$ws3.PageSetup.Orientation = 1
$ws3.Range("B2").Select() | Out-Null

$ws2.Activate() | Out-Null
$ws2.Range("B2").Select() | Out-Null
